$d = $word.ActiveDocument

# Each (find, replace) pair below corresponds to one run of text in the
# document: the date heading and the 25 division-problem table cells.
$pairs = @(
    @("2025-11-10 Monday", "2025-11-11 Tuesday"),
    @("301÷5=60, 1", "383÷9=42, 5"),
    @("822÷2=411, 0", "391÷4=97, 3"),
    @("662÷6=110, 2", "321÷6=53, 3"),
    @("186÷4=46, 2", "833÷6=138, 5"),
    @("376÷2=188, 0", "933÷4=233, 1"),
    @("999÷4=249, 3", "182÷6=30, 2"),
    @("400÷5=80, 0", "751÷3=250, 1"),
    @("905÷9=100, 5", "258÷2=129, 0"),
    @("514÷3=171, 1", "435÷3=145, 0"),
    @("552÷9=61, 3", "780÷8=97, 4"),
    @("157÷9=17, 4", "555÷8=69, 3"),
    @("755÷6=125, 5", "492÷7=70, 2"),
    @("245÷8=30, 5", "350÷9=38, 8"),
    @("968÷9=107, 5", "925÷7=132, 1"),
    @("963÷7=137, 4", "493÷2=246, 1"),
    @("547÷2=273, 1", "317÷7=45, 2"),
    @("259÷4=64, 3", "466÷3=155, 1"),
    @("166÷9=18, 4", "113÷4=28, 1"),
    @("321÷5=64, 1", "843÷2=421, 1"),
    @("997÷3=332, 1", "946÷3=315, 1"),
    @("649÷4=162, 1", "298÷7=42, 4"),
    @("465÷3=155, 0", "887÷8=110, 7"),
    @("689÷3=229, 2", "181÷8=22, 5"),
    @("938÷8=117, 2", "470÷9=52, 2"),
    @("398÷2=199, 0", "502÷8=62, 6"),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

Write-Output "Done."
